# Add team record (Wins/Losses/Ties) columns to the KCR_2009 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in row 1, columns AD, AE, AF.
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Copy the existing header formatting (bold / centered / bordered) from an
# existing header cell onto the new header cells so the same style is reused.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Team record is constant for every player row (2-45): 65 wins, 97 losses, 0 ties.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 65
    $ws.Cells.Item($r, 31).Value2 = 97
    $ws.Cells.Item($r, 32).Value2 = 0
}
